$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.380.42'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '1.567.27'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'1.001"

$ws.Range("D6").Value = "'291.01"
$ws.Range("E6").Value = '  +0.25%  '

$ws.Range("E7").Value = '  +1.99%  '

$ws.Range("D8").Value = "'49.05"
$ws.Range("E8").Value = '  -0.54%  '

$ws.Range("D10").Value = "'0.07568"
$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("D11").Value = "'1.134"
$ws.Range("E11").Value = '  -2.56%  '

$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = "'20.99"
$ws.Range("E13").Value = '  -1.66%  '

$ws.Range("E14").Value = '  -1.44%  '

$ws.Range("D15").Value = "'6.918"
$ws.Range("E15").Value = '  +0.14%  '

$ws.Range("D16").Value = '1.566.47'
$ws.Range("E16").Value = '  +0.10%  '

$ws.Range("E17").Value = '  +0.13%  '

$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("D19").Value = "'0.06750"
$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").Value = "'16.57"
$ws.Range("E21").Value = '  +0.43%  '

$ws.Range("D22").Value = "'6.196"
$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("E23").Value = '  -0.61%  '

$ws.Range("D24").Value = '22.368.40'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = "'2.381"
$ws.Range("E25").Value = '  +0.42%  '

$ws.Range("D26").Value = "'2.709"
$ws.Range("E26").Value = '  -3.90%  '

$ws.Range("D27").Value = "'20.17"
$ws.Range("E27").Value = '  +0.46%  '

$ws.Range("D28").Value = "'148.20"
$ws.Range("E28").Value = '  +1.13%  '

$ws.Range("D29").Value = "'5.034"
$ws.Range("E29").Value = '  +1.35%  '

$ws.Range("D30").Value = "'125.64"
$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").Value = '1.738.61'
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").Value = "'2.019"
$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'6.049"
$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'0.9893"
$ws.Range("E34").Value = '  -2.64%  '

$ws.Range("D35").Value = "'10.07"
$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("D36").Value = "'1.419"
$ws.Range("E36").Value = '  +11.39%  '

$ws.Range("D37").Value = "'0.08455"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("D38").Value = "'0.02485"
$ws.Range("E38").Value = '  -2.09%  '

$ws.Range("D39").Value = "'0.2288"
$ws.Range("E39").Value = '  -1.27%  '

$ws.Range("D40").Value = "'0.06453"
$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("D41").Value = "'5.414"
$ws.Range("E41").Value = '  -1.78%  '

$ws.Range("D42").Value = "'0.6301"
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("D43").Value = "'11.26"
$ws.Range("E43").Value = '  -3.42%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").Value = "'14.04"
$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("E46").Value = '  +1.17%  '

$ws.Range("D47").Value = "'0.5920"
$ws.Range("E47").Value = '  -0.84%  '

$ws.Range("D48").Value = "'2.072"
$ws.Range("E48").Value = '  -1.34%  '

$ws.Range("E49").Value = '  +0.17%  '

$ws.Range("D50").Value = "'124.66"
$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("D51").Value = "'0.07347"
$ws.Range("E51").Value = '  +0.88%  '
